$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 new rows at row 48, pushing current rows 48-101 down to 50-103.
$ws.Range("A48:R49").Insert()

# New row 48 data
$ws.Cells.Item(48, 1).Value = 11
$ws.Cells.Item(48, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(48, 3).Value = "Bíobío"
$ws.Cells.Item(48, 4).Value = 44923
$ws.Cells.Item(48, 5).Value = 8
$ws.Cells.Item(48, 6).Value = 100112012
$ws.Cells.Item(48, 7).Value = "Espinaca"
$ws.Cells.Item(48, 8).Value = "Sin especificar"
$ws.Cells.Item(48, 9).Value = "Primera"
$ws.Cells.Item(48, 10).Value = 100
$ws.Cells.Item(48, 11).Value = 7000
$ws.Cells.Item(48, 12).Value = 7500
$ws.Cells.Item(48, 13).Value = 7250
$ws.Cells.Item(48, 14).Value = "$/cuna 10 kilos"
$ws.Cells.Item(48, 15).Value = "Región Metropolitana"
$ws.Cells.Item(48, 16).Value = 725
$ws.Cells.Item(48, 17).Value = 10
$ws.Cells.Item(48, 18).Value = "Hortaliza"

# New row 49 data
$ws.Cells.Item(49, 1).Value = 11
$ws.Cells.Item(49, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(49, 3).Value = "Bíobío"
$ws.Cells.Item(49, 4).Value = 44923
$ws.Cells.Item(49, 5).Value = 8
$ws.Cells.Item(49, 6).Value = 100112012
$ws.Cells.Item(49, 7).Value = "Espinaca"
$ws.Cells.Item(49, 8).Value = "Sin especificar"
$ws.Cells.Item(49, 9).Value = "Primera"
$ws.Cells.Item(49, 10).Value = 100
$ws.Cells.Item(49, 11).Value = 7000
$ws.Cells.Item(49, 12).Value = 7500
$ws.Cells.Item(49, 13).Value = 7250
$ws.Cells.Item(49, 14).Value = "$/cuna 10 kilos"
$ws.Cells.Item(49, 15).Value = "Región Metropolitana"
$ws.Cells.Item(49, 16).Value = 725
$ws.Cells.Item(49, 17).Value = 10
$ws.Cells.Item(49, 18).Value = "Hortaliza"

# Make sure the date column (D) keeps the date number format for the new rows
$ws.Range("D48:D49").NumberFormat = "YYYY-MM-DD HH:MM:SS"
